$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume figures per latest GitHub Actions data pull.
# Each target cell holds its value as text (matches source inlineStr cells), so we
# force a Text number format before assigning to avoid Excel auto-converting the
# numeric-looking strings (prices, percentages) into actual numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '260.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.62%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.25'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.76%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.90%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06166'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.23%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.670'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.02%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8501'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.75%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9106'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.09%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.16%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04902'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '10.55%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07086'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.32%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03102'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.33%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09052'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.56%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.28%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006127'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.26%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006021'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.38%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.152'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.66%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.179'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.35%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.04%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.65%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.095'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.13%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04244'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.61%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.03%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003799'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-17.62%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.04%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001574'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-8.23%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03871'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.95%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.34%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004095'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-33.96%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '13.87%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002195'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.06%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005161'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.70%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '8.00%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1624'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-32.59%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
